# Update cryptocurrency price (D) and 1h volume change (E) columns
# with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.010.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.04%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.686.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.18%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.01'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.11%  '

$ws.Range("E10").Value = '  +0.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0884'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.78%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.926.94'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.685.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.12%  '

$ws.Range("E14").Value = '  -0.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.559'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.83'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '250.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.97%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.963.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0742'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.19%  '

$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("E24").Value = '  -2.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.23'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.17%  '

$ws.Range("E26").Value = '  -1.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.39%  '

$ws.Range("E28").Value = '  +0.24%  '

$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.27'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0503'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("E32").Value = '  +0.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.428.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.943'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.70%  '

$ws.Range("E37").Value = '  +0.27%  '

$ws.Range("E38").Value = '  -1.92%  '

$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("E40").Value = '  -2.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.42'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.37%  '

$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.833.58'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.797'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.27'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.72%  '

$ws.Range("E49").Value = '  -1.32%  '

$ws.Range("E50").Value = '  -1.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.52%  '
